$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PA73405 - Attrition by Job 2009")

# Grow the data table (Table3) from A1:F6 to A1:F8 to make room for the
# two new job rows ("Department Leader" and "Regional Leader").
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F8"))

# Columns A (year) and F (TerminationRate) hold values that look like a
# plain number ("2009") or a percentage ("30.8%") but must stay literal
# text, matching how the source data was authored. Force those columns to
# a text format before writing so Excel doesn't auto-convert them into a
# numeric/percentage value, then restore the Normal style so no visible
# formatting change is left behind.
$yearRange = $ws.Range("A2:A8")
$rateRange = $ws.Range("F2:F8")
$yearRange.NumberFormat = "@"
$rateRange.NumberFormat = "@"

$jobNames = @("Analyst", "Department Leader", "Regional Leader", "Product Manager", "Consultant", "Project Manager", "Administrative Assistant")
$didNotTerminate = @(9, 3, 11, 10, 22, 13, 16)
$terminated      = @(4, 1, 3, 2, 2, 1, 0)
$headcount       = @(13, 4, 14, 12, 24, 14, 16)
$terminationRate = @("30.8%", "25.0%", "21.4%", "16.7%", "8.3%", "7.1%", "0.0%")

for ($i = 0; $i -lt $jobNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "2009"
    $ws.Cells.Item($row, 2).Value = $jobNames[$i]
    $ws.Cells.Item($row, 3).Value = $didNotTerminate[$i]
    $ws.Cells.Item($row, 4).Value = $terminated[$i]
    $ws.Cells.Item($row, 5).Value = $headcount[$i]
    $ws.Cells.Item($row, 6).Value = $terminationRate[$i]
}

$yearRange.Style = "Normal"
$rateRange.Style = "Normal"
